# Add handling of multiple projects to td2cambio
#
# The "summa per projekt" row (row 39) on the "Konsulttidrapport" sheet only
# totalled the first project column (B). With multiple project columns
# (C..N) now in use, every project column needs its own monthly total, so
# fill C39:N39 with the same per-column SUM used in B39.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Konsulttidrapport")

# B39 already holds =SUM(B8:B38); extend the same relative formula across
# the rest of the project columns (C through N).
$ws.Range("C39:N39").FormulaR1C1 = "=SUM(R[-31]C:R[-1]C)"

# Reflect the author's updated selection/scroll position on the sheet.
$ws.Activate()
$ws.Range("B39:N39").Select()
